$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update existing amounts in column B
$ws.Range("B2").Value = 394
$ws.Range("B3").Value = 400

# Add new note in A12, same font as the row above (A11) plus a left-aligned date number format
$ws.Range("A12").Value = "19.02.2024- Otistics Portis karşılığında (ya da dize yatırma da diyebiliriz) Los Yahoo'ya 1 Dolar vermiştir. (394-400)"
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A12").HorizontalAlignment = -4131  # xlLeft
$ws.Range("A12").NumberFormat = "mm-dd-yy"

# Resize column A to fit the new (longer) content (best-fit width for the longer note)
$ws.Columns.Item(1).ColumnWidth = 87.3

# Restore selection to C5 as in the final saved view
$ws.Range("C5").Select()
